# Sync attendance_reports: normalize "Recorded By" (column G) values so that
# a leading "System, " prefix is moved to the end as ", System" instead.
#   e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#        "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$prefix = "System, "

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.StartsWith($prefix)) {
        $rest = $val.Substring($prefix.Length)
        $cell.Value = "$rest, System"
    }
}
